$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BASE DATA (wajib update)")

# The "Complete"/"Incomplete" helper columns (P and Q) previously carried a
# text placeholder (whitespace + "-") for rows where that metric didn't
# apply. The refreshed source data now reports these as a plain numeric 0,
# matching the rest of the column.
$cells = "P3,P5,P12,P13,P14,P15,P17,P22,P23,P24,P25,P26,P30,P31,P32,P36,P39,P48,P49,Q53,P54,Q55,Q57,Q58,Q62,P63,Q67,Q70,P73,P74,P77,P83,P86,Q89,P93,P95,P97,P101,P104,P107,P109,P110,P112,P116,Q117,P118,P121,P124"

$rng = $ws.Range($cells)
foreach ($area in $rng.Areas) {
    $area.Value = 0
}
